$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# "Given a station and range, find the best commodity to buy and sell
# locally." -> split the run so "station" becomes its own run, then swap
# its text for "system" (same run-splitting shape as the target diff: three
# runs "Given a " / "system" / " and range, ...").
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Given a station and range*") {
        $p1 = $p
        break
    }
}

$snippet1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="505AE2D4" w14:textId="77777777" w:rsidR="00E07FB3" w:rsidRPr="00E07FB3" w:rsidRDefault="00E07FB3" w:rsidP="00E07FB3">
<w:pPr>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
<w:spacing w:line="240" w:lineRule="auto"/>
<w:textAlignment w:val="baseline"/>
<w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Given a </w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>system</w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and range, find the best commodity to buy and sell locally.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p1.Range.InsertXML($snippet1)

# --- Edit 2 ---------------------------------------------------------------
# "Find the station to sell at that contains black market." -> "Find the
# nearest  station to sell at that contains black market." split into the
# five runs shown in the target diff: "Find the " / "nearest" / " " / " " /
# "station to sell at that contains black market."
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Find the station to sell at that contains black market*") {
        $p2 = $p
        break
    }
}

$snippet2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="68D8920C" w14:textId="77777777" w:rsidR="00E07FB3" w:rsidRPr="00E07FB3" w:rsidRDefault="00E07FB3" w:rsidP="00E07FB3">
<w:pPr>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>
<w:spacing w:line="240" w:lineRule="auto"/>
<w:textAlignment w:val="baseline"/>
<w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Find the </w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>nearest</w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>station to sell at that contains black market.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p2.Range.InsertXML($snippet2)
